# Swap the presentation's design theme palette from the "Integral" colour
# scheme over to the stock "Office Theme" colour scheme (this mirrors the
# underlying OOXML change where ppt/theme/theme1.xml <-> ppt/theme/theme2.xml
# content is swapped between the plain Office theme and the Integral theme).
#
# PowerPoint's object model doesn't give us raw part access, so we drive the
# same result through the Design/Theme/ThemeColorScheme COM surface that a
# user gets to via Design > Variants > Colors (Customize Colors...).

function To-RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Target palette = the standard Office Theme colours, in
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order (indices 1-12).
$officeColors = @(
    (To-RGB 0x00 0x00 0x00),   # 1  dk1
    (To-RGB 0xFF 0xFF 0xFF),   # 2  lt1
    (To-RGB 0x44 0x54 0x6A),   # 3  dk2
    (To-RGB 0xE7 0xE6 0xE6),   # 4  lt2
    (To-RGB 0x5B 0x9B 0xD5),   # 5  accent1
    (To-RGB 0xED 0x7D 0x31),   # 6  accent2
    (To-RGB 0xA5 0xA5 0xA5),   # 7  accent3
    (To-RGB 0xFF 0xC0 0x00),   # 8  accent4
    (To-RGB 0x44 0x72 0xC4),   # 9  accent5
    (To-RGB 0x70 0xAD 0x47),   # 10 accent6
    (To-RGB 0x05 0x63 0xC1),   # 11 hlink
    (To-RGB 0x95 0x4F 0x72)    # 12 folHlink
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
